$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scope (column C) updates on the Table1 rows ---

# Rows 23-27 ("cur.read" scope for the /v1/currencies* endpoints) no longer carry a Scope value
$ws.Range("C23:C27").ClearContents()

# Rows 35-42 previously had no Scope; they now read "ccs.read, ccs.write"
$ws.Range("C35:C42").Value = "ccs.read, ccs.write"

# Row 43 previously had no Scope; it now reads "ccs.read"
$ws.Range("C43").Value = "ccs.read"

# Rows 44-46 previously had no Scope; they now read "ccs.read, ccs.write"
$ws.Range("C44:C46").Value = "ccs.read, ccs.write"

# Row 49 narrows from "doc.read, doc.write" to "doc.read"
# (set before the cts.read rows below so new shared strings are introduced in
#  the same order as the authoritative edit: doc.read, then cts.read)
$ws.Range("C49").Value = "doc.read"

# Rows 28-29 narrow from "cts.read, cts.write" to "cts.read"
$ws.Range("C28:C29").Value = "cts.read"

# --- View / selection state ---
$window = $excel.ActiveWindow
$window.Zoom = 190
$ws.Range("A9").Select()
